# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the c90a607c... row, as part
# of regenerating the Handback report.
#
# These cells hold their date/time as plain text (they are formatted with a
# custom yyyy-mm-dd HH:mm:ss number format, but the underlying value is a
# text string), so write them back out as text too, to keep the same storage
# shape.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview!G4 - Latest HO Xliff Generate Date
$overview.Range("G4").Value = "2016-09-01 18:53:25"

# zh-cn!H4 - Correspond Handoff Datetime
$zhcn.Range("H4").Value = "2016-09-01 18:53:20"
# zh-cn!K4 - Correspond Handback DateTime
$zhcn.Range("K4").Value = "2016-09-01 18:53:57"

# de-de!H4 - Correspond Handoff Datetime (shares text w/ Overview!G4)
$dede.Range("H4").Value = "2016-09-01 18:53:25"
# de-de!K4 - Correspond Handback DateTime
$dede.Range("K4").Value = "2016-09-01 18:54:14"
